# Update the "Förändrad" (Changed) date column from 2023-10-06 (45205)
# to 2023-10-07 (45206) for every data row (rows 2 through 150).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C150").Value = 45206
